$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '99.289.99'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.370.82'
$ws.Range('E3').Value = '  +7.47%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.70'
$ws.Range('E5').Value = '  +7.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '629.27'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  +24.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.393'
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.864'
$ws.Range('E10').Value = '  +10.25%  '
$ws.Range('D11').Value = '3.368.60'
$ws.Range('E11').Value = '  +7.42%  '
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = '99.017.24'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.85'
$ws.Range('E14').Value = '  +4.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000247'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('D16').Value = '3.960.97'
$ws.Range('E16').Value = '  +6.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.50'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '3.367.74'
$ws.Range('E18').Value = '  +7.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.55'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.31'
$ws.Range('E20').Value = '  +4.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '494.86'
$ws.Range('E21').Value = '  -5.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.12'
$ws.Range('E22').Value = '  +7.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000209'
$ws.Range('E23').Value = '  +8.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.35'
$ws.Range('E24').Value = '  +5.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.66'
$ws.Range('E25').Value = '  +3.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.47'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.92'
$ws.Range('E27').Value = '  +2.56%  '
$ws.Range('D28').Value = '3.549.80'
$ws.Range('E28').Value = '  +7.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.279'
$ws.Range('E29').Value = '  +16.91%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.192'
$ws.Range('E31').Value = '  +9.72%  '
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +12.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.52'
$ws.Range('E34').Value = '  +5.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.98'
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.151'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.30'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '499.12'
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.458'
$ws.Range('E41').Value = '  +4.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.90'
$ws.Range('E42').Value = '  +8.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.26'
$ws.Range('E43').Value = '  +2.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.779'
$ws.Range('E46').Value = '  +10.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.04'
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('E49').Value = '  +13.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.63'
$ws.Range('E50').Value = '  +3.28%  '
$ws.Range('E51').Value = '  +4.16%  '
